$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace full team names with their standard 3-letter abbreviations.
# (ranks in column A and W-L% values in column C are unchanged)
$ws.Range("B6").Value  = "STL"
$ws.Range("B7").Value  = "NYY"
$ws.Range("B8").Value  = "BOS"
$ws.Range("B9").Value  = "ATL"
$ws.Range("B10").Value = "LAD"
$ws.Range("B11").Value = "MIN"
$ws.Range("B12").Value = "ANA"
$ws.Range("B13").Value = "HOU"
$ws.Range("B14").Value = "SFG"
$ws.Range("B15").Value = "OAK"
$ws.Range("B16").Value = "TEX"
$ws.Range("B17").Value = "CHC"
$ws.Range("B18").Value = "SDP"
$ws.Range("B19").Value = "PHI"
$ws.Range("B20").Value = "CWS"
$ws.Range("B21").Value = "FLA"
$ws.Range("B22").Value = "CLE"
$ws.Range("B23").Value = "BAL"
$ws.Range("B24").Value = "CIN"
$ws.Range("B25").Value = "PIT"
$ws.Range("B26").Value = "DET"
$ws.Range("B27").Value = "NYM"
$ws.Range("B28").Value = "TBD"
$ws.Range("B29").Value = "COL"
$ws.Range("B30").Value = "TOR"
$ws.Range("B31").Value = "MIL"
$ws.Range("B32").Value = "MON"
$ws.Range("B33").Value = "SEA"
$ws.Range("B34").Value = "KCR"
$ws.Range("B35").Value = "ARI"

# Stray blank marker cell that crept in next to row 8 (matches the " " shared string already used in A4)
$ws.Range("H8").Value = " "

# Sourcing footnote placed far below the table (row 43)
$ws.Range("H43").Value = 'Provided by <a href="https://www.sports-reference.com/sharing.html?utm_source=direct&utm_medium=Share&utm_campaign=ShareTool">Baseball-Reference.com</a>: <a href="https://www.baseball-reference.com/leagues/majors/2004-standings.shtml?sr&utm_source=direct&utm_medium=Share&utm_campaign=ShareTool#expanded_standings_overall">View Original Table</a><br>Generated 9/20/2023.'

# Update the view state: scroll/selection left where the user was last working
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("F11").Select()
